# Add a new row (11) of data to the "Word -> Replacement" lookup sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "hungvo"
$ws.Range("B11").Value = "123$5"
